$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44595
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100114002
$ws.Cells.Item($row, 7).Value = "Camote"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 50
$ws.Cells.Item($row, 11).Value = 18000
$ws.Cells.Item($row, 12).Value = 18000
$ws.Cells.Item($row, 13).Value = 18000
$ws.Cells.Item($row, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 900
$ws.Cells.Item($row, 17).Value = 20
$ws.Cells.Item($row, 18).Value = "Hortaliza"

$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
